$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9
$ws.Range("C2").Value = 'face/face079.png'
$ws.Range("D2").Value = 'kranken'
$ws.Range("E2").Value = 'face'

$ws.Range("B3").Value = 53
$ws.Range("C3").Value = 'face/face098.png'
$ws.Range("D3").Value = 'herrschen'
$ws.Range("E3").Value = 'face'

$ws.Range("B4").Value = 84
$ws.Range("C4").Value = 'face/face074.png'
$ws.Range("D4").Value = 'legen'
$ws.Range("E4").Value = 'face'

$ws.Range("B5").Value = 108
$ws.Range("C5").Value = 'face/face064.png'
$ws.Range("D5").Value = 'proben'
$ws.Range("E5").Value = 'face'

$ws.Range("B6").Value = 62
$ws.Range("C6").Value = 'flower/flower105.png'
$ws.Range("D6").Value = 'narren'
$ws.Range("E6").Value = 'flower'

$ws.Range("B7").Value = 86
$ws.Range("C7").Value = 'flower/flower087.png'
$ws.Range("D7").Value = 'binden'
$ws.Range("E7").Value = 'flower'

$ws.Range("B8").Value = 4
$ws.Range("C8").Value = 'face/face077.png'
$ws.Range("D8").Value = 'prüfen'
$ws.Range("E8").Value = 'face'

$ws.Range("B9").Value = 52
$ws.Range("C9").Value = 'flower/flower078.png'
$ws.Range("D9").Value = 'bergen'
$ws.Range("E9").Value = 'flower'

$ws.Range("B10").Value = 80
$ws.Range("C10").Value = 'face/face107.png'
$ws.Range("D10").Value = 'betteln'
$ws.Range("E10").Value = 'face'

$ws.Range("B11").Value = 123
$ws.Range("C11").Value = 'face/face116.png'
$ws.Range("D11").Value = 'kennen'
$ws.Range("E11").Value = 'face'

$ws.Range("B12").Value = 50
$ws.Range("C12").Value = 'face/face067.png'
$ws.Range("D12").Value = 'frischen'
$ws.Range("E12").Value = 'face'

$ws.Range("B13").Value = 13
$ws.Range("C13").Value = 'face/face082.png'
$ws.Range("D13").Value = 'fließen'
$ws.Range("E13").Value = 'face'

$ws.Range("B14").Value = 106
$ws.Range("C14").Value = 'flower/flower107.png'
$ws.Range("D14").Value = 'deuten'
$ws.Range("E14").Value = 'flower'

$ws.Range("B15").Value = 46
$ws.Range("C15").Value = 'face/face068.png'
$ws.Range("D15").Value = 'piepen'
$ws.Range("E15").Value = 'face'

$ws.Range("B16").Value = 111
$ws.Range("C16").Value = 'flower/flower072.png'
$ws.Range("D16").Value = 'zielen'
$ws.Range("E16").Value = 'flower'

$ws.Range("B17").Value = 22
$ws.Range("C17").Value = 'face/face072.png'
$ws.Range("D17").Value = 'leeren'
$ws.Range("E17").Value = 'face'

$ws.Range("B18").Value = 68
$ws.Range("C18").Value = 'flower/flower116.png'
$ws.Range("D18").Value = 'stoppen'
$ws.Range("E18").Value = 'flower'

$ws.Range("B19").Value = 45
$ws.Range("C19").Value = 'face/face106.png'
$ws.Range("D19").Value = 'nullen'
$ws.Range("E19").Value = 'face'

$ws.Range("B20").Value = 114
$ws.Range("C20").Value = 'face/face089.png'
$ws.Range("D20").Value = 'hören'
$ws.Range("E20").Value = 'face'

$ws.Range("B21").Value = 35
$ws.Range("C21").Value = 'face/face080.png'
$ws.Range("D21").Value = 'achten'
$ws.Range("E21").Value = 'face'

$ws.Range("B22").Value = 25
$ws.Range("C22").Value = 'flower/flower096.png'
$ws.Range("D22").Value = 'öffnen'
$ws.Range("E22").Value = 'flower'

$ws.Range("B23").Value = 40
$ws.Range("C23").Value = 'flower/flower088.png'
$ws.Range("D23").Value = 'planen'
$ws.Range("E23").Value = 'flower'

$ws.Range("B24").Value = 58
$ws.Range("C24").Value = 'flower/flower104.png'
$ws.Range("D24").Value = 'trotzen'
$ws.Range("E24").Value = 'flower'

$ws.Range("B25").Value = 66
$ws.Range("C25").Value = 'flower/flower094.png'
$ws.Range("D25").Value = 'reisen'
$ws.Range("E25").Value = 'flower'

$ws.Range("B26").Value = 33
$ws.Range("C26").Value = 'flower/flower086.png'
$ws.Range("D26").Value = 'lassen'
$ws.Range("E26").Value = 'flower'

$ws.Range("B27").Value = 105
$ws.Range("C27").Value = 'flower/flower103.png'
$ws.Range("D27").Value = 'küssen'
$ws.Range("E27").Value = 'flower'

$ws.Range("B28").Value = 83
$ws.Range("C28").Value = 'face/face075.png'
$ws.Range("D28").Value = 'rufen'
$ws.Range("E28").Value = 'face'

$ws.Range("B29").Value = 23
$ws.Range("C29").Value = 'flower/flower097.png'
$ws.Range("D29").Value = 'heben'
$ws.Range("E29").Value = 'flower'

$ws.Range("B30").Value = 51
$ws.Range("C30").Value = 'flower/flower083.png'
$ws.Range("D30").Value = 'lügen'
$ws.Range("E30").Value = 'flower'

$ws.Range("B31").Value = 37
$ws.Range("C31").Value = 'flower/flower082.png'
$ws.Range("D31").Value = 'ändern'
$ws.Range("E31").Value = 'flower'

$ws.Range("B32").Value = 54
$ws.Range("C32").Value = 'flower/flower064.png'
$ws.Range("D32").Value = 'spenden'
$ws.Range("E32").Value = 'flower'

$ws.Range("B33").Value = 127
$ws.Range("C33").Value = 'face/face099.png'
$ws.Range("D33").Value = 'mögen'
$ws.Range("E33").Value = 'face'
